$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Factor out a new "sc_isolation_entity list" sheet, positioned
#    right before "library_final_yield_unit list" (i.e. right after
#    "analyte_class list").
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "sc_isolation_entity list"

$libYieldSheet = $wb.Worksheets.Item("library_final_yield_unit list")
$newSheet.Move($libYieldSheet)

# Re-fetch a fresh reference post-move before writing values.
$entitySheet = $wb.Worksheets.Item("sc_isolation_entity list")
$entitySheet.Range("A1").Value = "whole cell"
$entitySheet.Range("A2").Value = "nucleus"
$entitySheet.Range("A3").Value = "cell-cell multimer"
$entitySheet.Range("A4").Value = "spatially encoded cell barcoding"

# ------------------------------------------------------------------
# 2. Add data validation on the sc_isolation_entity column (Q) of the
#    main "Export as TSV" sheet, pointing at the new list sheet.
# ------------------------------------------------------------------
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$val = $mainSheet.Range("Q2:Q1048576").Validation
$val.Add(3, 1, 1, "='sc_isolation_entity list'!`$A`$1:`$A`$4")
$val.ErrorTitle = "Value must come from list"
$val.ErrorMessage = "Value must be one of: whole cell / nucleus / cell-cell multimer / spatially encoded cell barcoding."
$val.ShowInput = $true
$val.ShowError = $true

# ------------------------------------------------------------------
# 3. Update the comment on column T (sc_isolation_quality_metric)
#    with the extended explanation.
# ------------------------------------------------------------------
$cmt = $mainSheet.Range("T1").Comment
$cmt.Text("A quality metric by visual inspection prior to cell lysis or defined by known parameters such as wells with several cells or no cells. This can be captured at a high level. `"OK`" or `"not OK`", or with more specificity such as `"debris`", `"clump`", `"low clump`".") | Out-Null
